# Apply "clean reports demo ready" changes to the
# protocoltestcasedetails sheet:
#   1. Add new column AJ "dataprofilelimit" with a value of 1000 for the
#      existing 30 test-case rows (2-31).
#   2. Flip the execute flags for testcase12 (row 13, C13 1 -> 0) and
#      testcase24 (row 25, C25 0 -> 1).
#   3. Fix the stage/target typo on row 17 (W17).
#   4. Append a brand new test case (row 32) for the snowflake parquet
#      validation scenario, including its AJ value of 2000.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("protocoltestcasedetails")

# --- 1. New "dataprofilelimit" column ------------------------------------
$ws.Range("AJ1").Value = "dataprofilelimit"

# Match the header formatting used by the rest of row 1 (bold, centered,
# bordered) by copying the style from the neighboring AI1 header cell.
$ws.Range("AI1").Copy()
$ws.Range("AJ1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AJ2:AJ31").Value = 1000

# --- 2. Toggle execute flags ----------------------------------------------
$ws.Range("C13").Value = 0
$ws.Range("C25").Value = 1

# --- 3. Correct target file path on row 17 --------------------------------
$ws.Range("W17").Value = "test/data/target"

# --- 4. Append new test case row 32 ---------------------------------------
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "testcase31_snowflake_parquet_validation"
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = "likeobjectcompare"
$ws.Range("E32").Value = "Auto"
$ws.Range("F32").Value = "rtpcr_source"
$ws.Range("G32").Value = "raw_snowflake_sql_connection"
$ws.Range("H32").Value = "snowflake"
$ws.Range("I32").Value = "table"
$ws.Range("K32").Value = "rtpcr_diagnostic_lab_testing"
$ws.Range("S32").Value = "rtpcr_target"
$ws.Range("U32").Value = "aws-s3"
$ws.Range("V32").Value = "parquet"
$ws.Range("W32").Value = "test/data/target"
$ws.Range("X32").Value = "patients_target_parquet_mismatch"
$ws.Range("AF32").Value = "test/s2t/s2t_31_snowflake_parquet_val.xlsx"
$ws.Range("AG32").Value = "source_to_target"
$ws.Range("AH32").Value = "id,state"
$ws.Range("AI32").Value = 8
$ws.Range("AJ32").Value = 2000
